$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Helper: wrap a <w:body> inner fragment in the minimal pkg:package /
# w:document shell that Range.InsertXML expects.
# ---------------------------------------------------------------------
function New-PkgXml([string]$bodyInner) {
    return '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' +
        "<w:body>$bodyInner</w:body>" +
        '</w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'
}

# Helper: find the first paragraph whose text contains $needle.
function Find-ParagraphContaining([string]$needle) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $cand = $d.Paragraphs.Item($i)
        if ($cand.Range.Text -like "*$needle*") { return $cand }
    }
    return $null
}

# Helper: grab the paragraph's own opening <w:p ...> tag (with its original
# w14:paraId / rsid attributes) straight from its current OOXML, so the
# replacement keeps those identifiers intact.
function Get-ParagraphOpenTag($para) {
    $xml = $para.Range.WordOpenXML
    if ($xml -match '(<w:p( [^>]*)?>)') {
        return $matches[1]
    }
    return "<w:p>"
}

# ---------------------------------------------------------------------
# Paragraph 1: the "{$ img: cad_model $}" placeholder paragraph collapses
# to a single empty run, preceded by a gramStart proof-error marker.
# ---------------------------------------------------------------------
$p1 = Find-ParagraphContaining "cad_model"
if ($p1 -ne $null) {
    $openTag1 = Get-ParagraphOpenTag $p1
    $inner1 = "$openTag1<w:proofErr w:type=`"gramStart`"/><w:r><w:t xml:space=`"preserve`"></w:t></w:r></w:p>"
    $p1.Range.InsertXML((New-PkgXml $inner1))
}

# ---------------------------------------------------------------------
# Paragraph 2: the "{$ img: plot $}" placeholder paragraph collapses to a
# single empty run (no proof-error marker).
# ---------------------------------------------------------------------
$p2 = Find-ParagraphContaining "img:plot"
if ($p2 -ne $null) {
    $openTag2 = Get-ParagraphOpenTag $p2
    $inner2 = "$openTag2<w:r><w:t xml:space=`"preserve`"></w:t></w:r></w:p>"
    $p2.Range.InsertXML((New-PkgXml $inner2))
}
